$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new values would otherwise be
# auto-converted to numbers by Excel (losing exact text representation).
$textCells = @("D5", "D8", "D10", "D11", "D13", "D14", "D15", "D19", "D21", "D26", "D27", "D28", "D31", "D32", "D33", "D36", "D37", "D38", "D40", "D47", "D48", "D50")
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Apply updated values
$ws.Range("D2").Value = "36.363.48"
$ws.Range("E2").Value = "  +0.19%  "
$ws.Range("D3").Value = "1.934.41"
$ws.Range("E3").Value = "  -2.28%  "
$ws.Range("D5").Value = "241.44"
$ws.Range("E5").Value = "  -1.70%  "
$ws.Range("E6").Value = "  -3.17%  "
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("D8").Value = "56.99"
$ws.Range("E8").Value = "  -3.87%  "
$ws.Range("E9").Value = "  -4.28%  "
$ws.Range("D10").Value = "0.0835"
$ws.Range("E10").Value = "  +1.33%  "
$ws.Range("D11").Value = "0.104"
$ws.Range("E11").Value = "  -0.31%  "
$ws.Range("D12").Value = "2.217.65"
$ws.Range("E12").Value = "  -2.27%  "
$ws.Range("D13").Value = "0.802"
$ws.Range("E13").Value = "  -6.91%  "
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").Value = "13.41"
$ws.Range("E14").Value = "  -3.75%  "
$ws.Range("B15").Value = "Avalanche"
$ws.Range("C15").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D15").Value = "20.97"
$ws.Range("E15").Value = "  -10.97%  "
$ws.Range("E16").Value = "  -5.64%  "
$ws.Range("D17").Value = "1.948.52"
$ws.Range("E17").Value = "  -1.61%  "
$ws.Range("D18").Value = "36.267.70"
$ws.Range("E18").Value = "  +0.06%  "
$ws.Range("D19").Value = "68.94"
$ws.Range("E19").Value = "  -1.46%  "
$ws.Range("D20").Value = "0.0₃0862"
$ws.Range("E20").Value = "  -0.66%  "
$ws.Range("D21").Value = "227.39"
$ws.Range("E21").Value = "  -2.91%  "
$ws.Range("E22").Value = "  -6.66%  "
$ws.Range("E23").Value = "  -0.07%  "
$ws.Range("E24").Value = "  -10.68%  "
$ws.Range("E25").Value = "  -1.74%  "
$ws.Range("D26").Value = "9.28"
$ws.Range("E26").Value = "  -7.28%  "
$ws.Range("D27").Value = "160.79"
$ws.Range("E27").Value = "  -0.75%  "
$ws.Range("D28").Value = "0.130"
$ws.Range("E28").Value = "  -1.46%  "
$ws.Range("E29").Value = "  -3.09%  "
$ws.Range("E30").Value = "  -2.43%  "
$ws.Range("D31").Value = "1.12"
$ws.Range("E31").Value = "  -5.91%  "
$ws.Range("D32").Value = "4.56"
$ws.Range("E32").Value = "  -6.81%  "
$ws.Range("D33").Value = "0.0626"
$ws.Range("E33").Value = "  -0.03%  "
$ws.Range("E34").Value = "  -5.54%  "
$ws.Range("E35").Value = "  +0.01%  "
$ws.Range("D36").Value = "6.07"
$ws.Range("E36").Value = "  -2.38%  "
$ws.Range("D37").Value = "1.79"
$ws.Range("E37").Value = "  -0.36%  "
$ws.Range("D38").Value = "2.13"
$ws.Range("E38").Value = "  -5.94%  "
$ws.Range("E39").Value = "  -1.45%  "
$ws.Range("D40").Value = "0.0970"
$ws.Range("E40").Value = "  +0.45%  "
$ws.Range("E41").Value = "  -1.13%  "
$ws.Range("E42").Value = "  -7.39%  "
$ws.Range("E43").Value = "  -2.75%  "
$ws.Range("E44").Value = "  -4.10%  "
$ws.Range("D45").Value = "1.331.96"
$ws.Range("E45").Value = "  -2.69%  "
$ws.Range("E46").Value = "  -7.19%  "
$ws.Range("D47").Value = "86.13"
$ws.Range("E47").Value = "  -6.67%  "
$ws.Range("D48").Value = "7.07"
$ws.Range("E48").Value = "  -5.53%  "
$ws.Range("E49").Value = "  -0.69%  "
$ws.Range("D50").Value = "44.05"
$ws.Range("E50").Value = "  -3.11%  "
$ws.Range("D51").Value = "2.108.78"
$ws.Range("E51").Value = "  -2.39%  "
